# Auto-generated: apply scheduled market-data refresh to Leve profit sheets
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H-N) per row
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 13161531
$ws.Range("I33").Value = 19231384
$ws.Range("J33").Value = 10185
$ws.Range("K33").Value = 19231384
$ws.Range("L33").Value = 10185
$ws.Range("M33").Value = -19231155
$ws.Range("N33").Value = -10643
$ws.Range("H69").Value = 26503.5
$ws.Range("I69").Value = 3000
$ws.Range("J69").Value = 50007
$ws.Range("K69").Value = 9000
$ws.Range("L69").Value = 150021
$ws.Range("M69").Value = -8126
$ws.Range("N69").Value = -151769
$ws.Range("H72").Value = 26503.5
$ws.Range("I72").Value = 3000
$ws.Range("J72").Value = 50007
$ws.Range("K72").Value = 27000
$ws.Range("L72").Value = 450063
$ws.Range("M72").Value = -22632
$ws.Range("N72").Value = -458799
$ws.Range("H107").Value = 1422
$ws.Range("J107").Value = 1974
$ws.Range("L107").Value = 1974
$ws.Range("N107").Value = -5814
$ws.Range("H111").Value = 846.125
$ws.Range("I111").Value = 552.3333
$ws.Range("J111").Value = 1727.5
$ws.Range("K111").Value = 1656.9999
$ws.Range("L111").Value = 5182.5
$ws.Range("M111").Value = 1410.0001
$ws.Range("N111").Value = -11316.5

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 10047.105
$ws.Range("I32").Value = 3425.1965
$ws.Range("J32").Value = 28588.45
$ws.Range("K32").Value = 3425.1965
$ws.Range("L32").Value = 28588.45
$ws.Range("M32").Value = -3138.1965
$ws.Range("N32").Value = -29162.45

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 2721.4211
$ws.Range("I86").Value = 2247.1333
$ws.Range("J86").Value = 4500
$ws.Range("K86").Value = 2247.1333
$ws.Range("L86").Value = 4500
$ws.Range("M86").Value = -1124.1333
$ws.Range("N86").Value = -6746
$ws.Range("H89").Value = 2721.4211
$ws.Range("I89").Value = 2247.1333
$ws.Range("J89").Value = 4500
$ws.Range("K89").Value = 11235.6665
$ws.Range("L89").Value = 22500
$ws.Range("M89").Value = -5619.666499999999
$ws.Range("N89").Value = -33732
$ws.Range("H94").Value = 4555.273
$ws.Range("I94").Value = 2036.8572
$ws.Range("K94").Value = 2036.8572
$ws.Range("M94").Value = -1585.8572
$ws.Range("H134").Value = 12280.906
$ws.Range("I134").Value = 6320.0356
$ws.Range("K134").Value = 18960.1068
$ws.Range("M134").Value = -16425.1068

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 28132.87
$ws.Range("I31").Value = 18955.166
$ws.Range("K31").Value = 18955.166
$ws.Range("M31").Value = -18660.166
$ws.Range("H34").Value = 28132.87
$ws.Range("I34").Value = 18955.166
$ws.Range("K34").Value = 18955.166
$ws.Range("M34").Value = -18753.166
$ws.Range("H62").Value = 6508.2
$ws.Range("J62").Value = 6241.5
$ws.Range("L62").Value = 6241.5
$ws.Range("N62").Value = -7489.5
$ws.Range("H65").Value = 6508.2
$ws.Range("J65").Value = 6241.5
$ws.Range("L65").Value = 31207.5
$ws.Range("N65").Value = -37447.5
$ws.Range("H74").Value = 41666.668
$ws.Range("J74").Value = 41666.668
$ws.Range("L74").Value = 41666.668
$ws.Range("N74").Value = -43414.668
$ws.Range("H77").Value = 41666.668
$ws.Range("J77").Value = 41666.668
$ws.Range("L77").Value = 125000.004
$ws.Range("N77").Value = -133736.004

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 381.11765
$ws.Range("I2").Value = 433.83334
$ws.Range("K2").Value = 2603.00004
$ws.Range("M2").Value = -2490.00004
$ws.Range("H12").Value = 37.42857
$ws.Range("I12").Value = 0
$ws.Range("J12").Value = 37.42857
$ws.Range("K12").Value = 0
$ws.Range("L12").Value = 112.28571
$ws.Range("M12").ClearContents()
$ws.Range("N12").Value = -458.28571
$ws.Range("H50").Value = 1596235.8
$ws.Range("I50").Value = 2552
$ws.Range("J50").Value = 2506912.2
$ws.Range("K50").Value = 7656
$ws.Range("L50").Value = 7520736.600000001
$ws.Range("M50").Value = -7175
$ws.Range("N50").Value = -7521698.600000001
$ws.Range("H53").Value = 1596235.8
$ws.Range("I53").Value = 2552
$ws.Range("J53").Value = 2506912.2
$ws.Range("K53").Value = 7656
$ws.Range("L53").Value = 7520736.600000001
$ws.Range("M53").Value = -7175
$ws.Range("N53").Value = -7521698.600000001
$ws.Range("H121").Value = 2356.2856
$ws.Range("J121").Value = 3498.5
$ws.Range("L121").Value = 10495.5
$ws.Range("N121").Value = -13115.5
$ws.Range("H128").Value = 184970.75
$ws.Range("I128").Value = 184970.75
$ws.Range("K128").Value = 554912.25
$ws.Range("M128").Value = -549932.25
$ws.Range("H131").Value = 1485.01
$ws.Range("J131").Value = 1492.6364
$ws.Range("L131").Value = 4477.9092
$ws.Range("N131").Value = -14557.9092

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H5").Value = 3208
$ws.Range("I5").Value = 3208
$ws.Range("J5").Value = 0
$ws.Range("K5").Value = 3208
$ws.Range("L5").Value = 0
$ws.Range("M5").Value = -3096
$ws.Range("N5").ClearContents()
$ws.Range("H102").Value = 9087.799999999999
$ws.Range("I102").Value = 6178.3076
$ws.Range("K102").Value = 6178.3076
$ws.Range("M102").Value = -4556.3076

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 4042.52
$ws.Range("I61").Value = 2718.182
$ws.Range("J61").Value = 5083.0713
$ws.Range("K61").Value = 2718.182
$ws.Range("L61").Value = 5083.0713
$ws.Range("M61").Value = -2516.182
$ws.Range("N61").Value = -5487.0713
$ws.Range("H68").Value = 3766496
$ws.Range("I68").Value = 496.66666
$ws.Range("J68").Value = 4472621
$ws.Range("K68").Value = 496.66666
$ws.Range("L68").Value = 4472621
$ws.Range("M68").Value = 252.33334
$ws.Range("N68").Value = -4474119
$ws.Range("H71").Value = 3766496
$ws.Range("I71").Value = 496.66666
$ws.Range("J71").Value = 4472621
$ws.Range("K71").Value = 2483.3333
$ws.Range("L71").Value = 22363105
$ws.Range("M71").Value = 1260.6667
$ws.Range("N71").Value = -22370593
$ws.Range("H107").Value = 9266.538
$ws.Range("I107").Value = 9266.538
$ws.Range("K107").Value = 9266.538
$ws.Range("M107").Value = -7346.538
$ws.Range("H113").Value = 4042.52
$ws.Range("I113").Value = 2718.182
$ws.Range("J113").Value = 5083.0713
$ws.Range("K113").Value = 2718.182
$ws.Range("L113").Value = 5083.0713
$ws.Range("M113").Value = -548.1819999999998
$ws.Range("N113").Value = -9423.0713

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 3643.1177
$ws.Range("I62").Value = 3386.0833
$ws.Range("J62").Value = 4260
$ws.Range("K62").Value = 3386.0833
$ws.Range("L62").Value = 4260
$ws.Range("M62").Value = -2762.0833
$ws.Range("N62").Value = -5508
$ws.Range("H65").Value = 3643.1177
$ws.Range("I65").Value = 3386.0833
$ws.Range("J65").Value = 4260
$ws.Range("K65").Value = 16930.4165
$ws.Range("L65").Value = 21300
$ws.Range("M65").Value = -13810.4165
$ws.Range("N65").Value = -27540
$ws.Range("H98").Value = 36266.668
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 36266.668
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 36266.668
$ws.Range("M98").ClearContents()
$ws.Range("N98").Value = -42256.668
$ws.Range("H126").Value = 8479.5
$ws.Range("I126").Value = 5392.9165
$ws.Range("K126").Value = 16178.7495
$ws.Range("M126").Value = -13708.7495
